# Append the new quarterly data point (row 84) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the date-cell formatting used by the rest of column A (e.g. A83) -
# copy formats only so no new style entries are introduced.
$ws.Range("A83").Copy()
$ws.Range("A84").PasteSpecial(-4122)

# New row: date 2025-08-15 (serial 45884) in column A, value -0.3 in column B.
$ws.Range("A84").Value = 45884
$ws.Range("B84").Value = -0.3
